$d = $word.ActiveDocument

$replacements = @(
    @{old = "750÷9=83, 3"; new = "575÷2=287, 1"},
    @{old = "825÷5=165, 0"; new = "174÷6=29, 0"},
    @{old = "173÷7=24, 5"; new = "841÷7=120, 1"},
    @{old = "227÷4=56, 3"; new = "513÷4=128, 1"},
    @{old = "646÷8=80, 6"; new = "310÷2=155, 0"},
    @{old = "435÷8=54, 3"; new = "632÷9=70, 2"},
    @{old = "272÷2=136, 0"; new = "640÷4=160, 0"},
    @{old = "230÷3=76, 2"; new = "151÷8=18, 7"},
    @{old = "464÷8=58, 0"; new = "100÷4=25, 0"},
    @{old = "718÷2=359, 0"; new = "983÷7=140, 3"},
    @{old = "407÷3=135, 2"; new = "270÷5=54, 0"},
    @{old = "315÷8=39, 3"; new = "731÷7=104, 3"},
    @{old = "603÷6=100, 3"; new = "731÷2=365, 1"},
    @{old = "730÷8=91, 2"; new = "247÷5=49, 2"},
    @{old = "237÷5=47, 2"; new = "779÷2=389, 1"},
    @{old = "838÷9=93, 1"; new = "454÷6=75, 4"},
    @{old = "814÷7=116, 2"; new = "629÷9=69, 8"},
    @{old = "552÷5=110, 2"; new = "854÷8=106, 6"},
    @{old = "159÷8=19, 7"; new = "747÷4=186, 3"},
    @{old = "906÷7=129, 3"; new = "636÷8=79, 4"},
    @{old = "120÷6=20, 0"; new = "980÷6=163, 2"},
    @{old = "460÷2=230, 0"; new = "600÷2=300, 0"},
    @{old = "817÷5=163, 2"; new = "251÷5=50, 1"},
    @{old = "830÷6=138, 2"; new = "900÷6=150, 0"},
    @{old = "273÷9=30, 3"; new = "630÷4=157, 2"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
